# Auto-generated PowerShell COM-interop script
# Applies "Update gh-pages to output generated at 456a3b4" diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 200
$ws.Range("F3").Value = 2468
$ws.Range("F5").Value = 1799
$ws.Range("F6").Value = 113
$ws.Range("F7").Value = 327
$ws.Range("F9").Value = 3596
$ws.Range("F10").Value = 1214
$ws.Range("F12").Value = 31
$ws.Range("F15").Value = 1396
$ws.Range("F17").Value = 1801
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 473
$ws.Range("F21").Value = 1561
$ws.Range("F23").Value = 7
$ws.Range("F25").Value = 2329
$ws.Range("F26").Value = 268
$ws.Range("F28").Value = 4372
$ws.Range("F35").Value = 949
$ws.Range("F36").Value = 10

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F22").Value = 149
$ws.Range("F29").Value = 67
$ws.Range("F35").Value = 450
$ws.Range("F42").Value = 89
$ws.Range("F45").Value = 74
$ws.Range("F46").Value = 30
$ws.Range("F47").Value = 30

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2546
$ws.Range("F4").Value = 2557
$ws.Range("F5").Value = 9586
$ws.Range("F9").Value = 396
$ws.Range("F10").Value = 3022
$ws.Range("F11").Value = 530
$ws.Range("F13").Value = 262

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")

# Row 2: interest count bump only
$ws.Range("F2").Value = 2546

# Shift rows 3-10 down to 4-11 (a new row is inserted at 3;
# the old row 11 "爆裂鼓手" event drops out of this curated sheet)
# old row 10 content -> row 11
$ws.Range("B11").Value = "2024-10-17"
$ws.Range("C11").Value = "上海·蜡笔小新：我们的恐龙日记x HAPPY ZOO 主题咖啡厅"
$ws.Range("D11").Value = "南京东路340号百联zx创趣场四楼05号 HAPPY ZOO"
$ws.Range("E11").Value = "2024.10.17 00:00-10.27 23:59"
$ws.Range("F11").Value = 258
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=93221"
$ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202410/nzGP5KRA1728526131597.png"

# old row 9 content -> row 10
$ws.Range("B10").Value = "2024-10-10"
$ws.Range("C10").Value = "上海·「火影忍者疾风传 × animate cafe」"
$ws.Range("D10").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws.Range("E10").Value = "2024.10.10 00:00-11.12 23:59"
$ws.Range("F10").Value = 823
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=92883"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202409/aQIhaIgt1727249498713.png"

# old row 8 content -> row 9
$ws.Range("B9").Value = "2024-10-01"
$ws.Range("C9").Value = "上海·2024·《世界之外》x  萌果酱谷子咖啡"
$ws.Range("D9").Value = "南京东路340号百联ZX 萌果酱谷子咖啡（百联）"
$ws.Range("E9").Value = "2024.10.01 00:00-12.11 23:59"
$ws.Range("F9").Value = 3020
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=93006"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202409/qtffZOKB1727426243733.png"

# old row 7 content -> row 8
$ws.Range("B8").Value = "2024-09-26"
$ws.Range("C8").Value = "上海·【神秘的西夏陵】大空间高沉浸探险体验"
$ws.Range("D8").Value = "南京西路325号 上海市历史博物馆"
$ws.Range("E8").Value = "2024.09.26 10:00-12.31 19:00"
$ws.Range("F8").Value = 33
$ws.Range("G8").Value = 108
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=92581"
$ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202409/jeDZO2cS1726302714881.jpeg"

# old row 6 content -> row 7
$ws.Range("B7").Value = "2024-09-24"
$ws.Range("C7").Value = "上海·星零界·社交游乐·休闲运动·潮玩派对"
$ws.Range("D7").Value = "长宁路1191号长宁来福士B1 长宁来福士"
$ws.Range("E7").Value = "2024.09.24 10:00-12.31 22:00"
$ws.Range("F7").Value = 16
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=92659"
$ws.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202409/PHS8s1lu1726221065737.png"

# old row 5 content -> row 6
$ws.Range("B6").Value = "2024-09-15"
$ws.Range("C6").Value = "上海 洛天依歌行宇宙·无限遨游 沉浸式体验展"
$ws.Range("D6").Value = "中山北路3300号 上海月星环球港"
$ws.Range("E6").Value = "2024.09.15 10:00-10.31 22:00"
$ws.Range("F6").Value = 2467
$ws.Range("G6").Value = 138
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=91175"
$ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202408/ei9COXS41724405861343.jpeg"

# old row 4 content -> row 5
$ws.Range("B5").Value = "2024-09-14"
$ws.Range("C5").Value = "上海·吉卜力工作室物语-沉浸式艺术展全球首站（9月-10月）"
$ws.Range("D5").Value = "龙台路10号2F 上海国际传媒港艺术中心"
$ws.Range("E5").Value = "2024.09.14 10:00-10.31 20:00"
$ws.Range("F5").Value = 198
$ws.Range("G5").Value = 9.9
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=91856"
$ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202409/wSR0yFfg1725432304586.jpeg"

# old row 3 content -> row 4
$ws.Range("B4").Value = "2024-09-10"
$ws.Range("C4").Value = "上海·迷你四驱车赛场"
$ws.Range("D4").Value = "虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶"
$ws.Range("E4").Value = "2024.09.10 10:00-12.31 22:00"
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 48
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=92042"
$ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202409/LzFT5TMO1725348229429.png"

# New row 3: 上海·东方明珠 event (newly surfaced in this curated sheet)
$ws.Range("B3").Value = "2024-08-17"
$ws.Range("C3").Value = "上海·东方明珠·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题店"
$ws.Range("D3").Value = "世纪大道1号 东方明珠电视塔城市广场商场"
$ws.Range("E3").Value = "2024.08.17 00:00-10.27 23:59"
$ws.Range("F3").Value = 2557
$ws.Range("G3").Value = "已售罄"
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=90444"
$ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202408/qUE9n4UR1723020534077.png"

# Rows 12-51: content unchanged, interest-count values updated
$ws.Range("F12").Value = 1799
$ws.Range("F13").Value = 113
$ws.Range("F14").Value = 327
$ws.Range("F16").Value = 3596
$ws.Range("F18").Value = 1214
$ws.Range("F28").Value = 1801
$ws.Range("F29").Value = 1561
$ws.Range("F31").Value = 149
$ws.Range("F32").Value = 149
$ws.Range("F33").Value = 7
$ws.Range("F39").Value = 4372
$ws.Range("F41").Value = 450
$ws.Range("F48").Value = 74
$ws.Range("F49").Value = 30
$ws.Range("F50").Value = 949
$ws.Range("F51").Value = 10

